# Update cryptocurrency price/volume data on Sheet1 (rows 2-51).
# Mirrors the "Updated cryptos list" GitHub Actions scrape commit:
# price (D) and 1h volume (E) values refreshed; a few rows also swapped
# rank position (name/link/price/volume all updated in place) to reflect
# the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.178.55"
$ws.Range("E2").Value = "  +10.98%  "
$ws.Range("D3").Value = "1.677.28"
$ws.Range("E3").Value = "  +6.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.87"
$ws.Range("E5").Value = "  +7.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3728"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3444"
$ws.Range("E8").Value = "  +5.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.45"
$ws.Range("E9").Value = "  +13.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.186"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07305"
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.47"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.096"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.769"
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("D16").Value = "1.682.31"
$ws.Range("E16").Value = "  +7.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06709"
$ws.Range("E19").Value = "  +7.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.71"
$ws.Range("E20").Value = "  +9.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.45"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.096"
$ws.Range("E22").Value = "  +3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.01"
$ws.Range("E23").Value = "  +3.34%  "
$ws.Range("D24").Value = "24.161.74"
$ws.Range("E24").Value = "  +10.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.367"
$ws.Range("E26").Value = "  -9.15%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.673"
$ws.Range("E27").Value = "  +12.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.23"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.57"
$ws.Range("E29").Value = "  +6.12%  "
$ws.Range("D30").Value = "1.869.29"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.85"
$ws.Range("E31").Value = "  +5.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.408"
$ws.Range("E32").Value = "  +17.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.109"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9817"
$ws.Range("E34").Value = "  +7.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.782"
$ws.Range("E35").Value = "  +11.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08458"
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.34"
$ws.Range("E37").Value = "  +6.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06457"
$ws.Range("E38").Value = "  +6.34%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.356"
$ws.Range("E39").Value = "  +4.28%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.903"
$ws.Range("E40").Value = "  +7.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02338"
$ws.Range("E41").Value = "  +7.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.272"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2122"
$ws.Range("E43").Value = "  +5.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6172"
$ws.Range("E44").Value = "  +7.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9987"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.800"
$ws.Range("E46").Value = "  +4.76%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.18"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5949"
$ws.Range("E48").Value = "  +5.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.82"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.032"
$ws.Range("E50").Value = "  +6.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07183"
$ws.Range("E51").Value = "  +6.38%  "
